$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Spatial Impasse"
$ws.Range("I8").Value = 1

$ws.Range("I11").Select()
